$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DC cable existing capacity matrix values
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 2
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("I5").Value = 0

$ws.Range("D6").Value = 0
$ws.Range("G6").Value = 4

$ws.Range("F7").Value = 4

$ws.Range("B8").Value = 4

$ws.Range("B9").Value = 0
$ws.Range("E9").Value = 0

# Match the active selection recorded in the saved workbook
$ws.Range("F8").Select()
